$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.859.70"
$ws.Range("E2").Value = "  -0.14%  "

# Row 3
$ws.Range("D3").Value = "'1.876.64"
$ws.Range("E3").Value = "  +0.15%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'0.7263"
$ws.Range("E5").Value = "  -1.56%  "

# Row 6
$ws.Range("D6").Value = "'241.77"
$ws.Range("E6").Value = "  -0.24%  "

# Row 7
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").Value = "'0.3151"
$ws.Range("E8").Value = "  -0.11%  "

# Row 9
$ws.Range("D9").Value = "'0.07681"
$ws.Range("E9").Value = "  +7.02%  "

# Row 10
$ws.Range("D10").Value = "'24.62"
$ws.Range("E10").Value = "  -0.18%  "

# Row 11
$ws.Range("D11").Value = "'0.08174"
$ws.Range("E11").Value = "  -2.85%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.929.49"
$ws.Range("E12").Value = "  +2.95%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7469"
$ws.Range("E13").Value = "  -0.43%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.357"
$ws.Range("E14").Value = "  -1.13%  "

# Row 15
$ws.Range("D15").Value = "'92.70"
$ws.Range("E15").Value = "  +0.19%  "

# Row 16
$ws.Range("D16").Value = "'29.916.30"
$ws.Range("E16").Value = "  +0.03%  "

# Row 17
$ws.Range("D17").Value = "'6.029"
$ws.Range("E17").Value = "  -1.13%  "

# Row 18
$ws.Range("D18").Value = "'247.58"
$ws.Range("E18").Value = "  +2.00%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007996"
$ws.Range("E19").Value = "  +2.36%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'13.47"
$ws.Range("E20").Value = "  -0.77%  "

# Row 21
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.23%  "

# Row 22
$ws.Range("D22").Value = "'2.140.99"
$ws.Range("E22").Value = "  +0.80%  "

# Row 23
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.05%  "

# Row 24
$ws.Range("D24").Value = "'7.751"
$ws.Range("E24").Value = "  -3.07%  "

# Row 25
$ws.Range("D25").Value = "'0.1522"
$ws.Range("E25").Value = "  -1.81%  "

# Row 26
$ws.Range("D26").Value = "'9.273"
$ws.Range("E26").Value = "  -0.12%  "

# Row 27
$ws.Range("D27").Value = "'164.19"
$ws.Range("E27").Value = "  -0.44%  "

# Row 28
$ws.Range("E28").Value = "  +0.18%  "

# Row 29
$ws.Range("D29").Value = "'2.012"
$ws.Range("E29").Value = "  -1.08%  "

# Row 30
$ws.Range("D30").Value = "'1.439"
$ws.Range("E30").Value = "  -3.44%  "

# Row 31
$ws.Range("D31").Value = "'4.540"
$ws.Range("E31").Value = "  -1.21%  "

# Row 32
$ws.Range("E32").Value = "  -0.05%  "

# Row 33
$ws.Range("D33").Value = "'4.202"
$ws.Range("E33").Value = "  -1.09%  "

# Row 34
$ws.Range("D34").Value = "'0.05442"
$ws.Range("E34").Value = "  +2.34%  "

# Row 35
$ws.Range("D35").Value = "'1.232"
$ws.Range("E35").Value = "  -0.21%  "

# Row 36
$ws.Range("D36").Value = "'0.7428"
$ws.Range("E36").Value = "  -1.51%  "

# Row 37
$ws.Range("E37").Value = "  +0.51%  "

# Row 38
$ws.Range("D38").Value = "'2.699"
$ws.Range("E38").Value = "  +0.27%  "

# Row 39
$ws.Range("D39").Value = "'0.01930"
$ws.Range("E39").Value = "  -0.88%  "

# Row 40
$ws.Range("D40").Value = "'2.743"
$ws.Range("E40").Value = "  -0.39%  "

# Row 41
$ws.Range("D41").Value = "'0.4476"
$ws.Range("E41").Value = "  -0.77%  "

# Row 42
$ws.Range("D42").Value = "'0.8863"
$ws.Range("E42").Value = "  +3.39%  "

# Row 43
$ws.Range("D43").Value = "'5.992"
$ws.Range("E43").Value = "  -1.04%  "

# Row 44
$ws.Range("D44").Value = "'71.86"
$ws.Range("E44").Value = "  -0.45%  "

# Row 45
$ws.Range("D45").Value = "'104.37"
$ws.Range("E45").Value = "  +1.16%  "

# Row 46
$ws.Range("D46").Value = "'1.043.14"
$ws.Range("E46").Value = "  -6.04%  "

# Row 47
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  -0.10%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.736"
$ws.Range("E48").Value = "  +2.12%  "

# Row 49
$ws.Range("D49").Value = "'1.825"
$ws.Range("E49").Value = "  -0.76%  "

# Row 50
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'7.488"
$ws.Range("E50").Value = "  -2.10%  "

# Row 51
$ws.Range("D51").Value = "'2.026.73"
$ws.Range("E51").Value = "  +0.53%  "
